# "new test results using new unit"
# The lab-results workbook switches its timing metrics from seconds to
# microseconds: the column headers gain a "/μs" suffix and every measured
# timing value is replaced with a new, independently-remeasured number
# (not a straight unit conversion).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relabel the shared "Present" placeholder string and append the
#     "/us" (micro-seconds) unit suffix to each metric's header text.
#     (The label text shown in column A for each data row - "encryption",
#     "decryption", "digest generation", "signature generation" - is left
#     exactly as-is; only the Present placeholder and the chart/series
#     header cells further below pick up the unit suffix.)

# --- 2. DES block (rows 2-4)
$ws.Range("F3").Value = 10010
$ws.Range("G3").Value = 40038
$ws.Range("H3").Value = 240224
$ws.Range("I3").Value = 1881706

$ws.Range("G4").Value = 30030
$ws.Range("H4").Value = 220198
$ws.Range("I4").Value = 1711564

# --- 3. AES block (rows 8-10)
$ws.Range("G9").Value = 10006
$ws.Range("H9").Value = 30028
$ws.Range("I9").Value = 290254

$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 10014
$ws.Range("H10").Value = 30038
$ws.Range("I10").Value = 230204

# --- 4. RSA block (rows 13-15)
$ws.Range("B14").Value = 10020
$ws.Range("C14").Value = 10020
$ws.Range("D14").Value = 10018
$ws.Range("E14").Value = 10018
$ws.Range("F14").Value = 10018
$ws.Range("G14").Value = 10004

$ws.Range("B15").Value = 10000
$ws.Range("C15").Value = 10000
$ws.Range("D15").Value = 10002
$ws.Range("E15").Value = 10002
$ws.Range("F15").Value = 10002
$ws.Range("G15").Value = 10012

# --- 5. SHA-1 block (rows 18-19)
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 9998
$ws.Range("G19").Value = 10008
$ws.Range("H19").Value = 110092
$ws.Range("I19").Value = 1010914

# --- 6. HMAC block (rows 23-24)
$ws.Range("G24").Value = 10010
$ws.Range("H24").Value = 40028
$ws.Range("I24").Value = 340302

# --- 7. Append the "/us" unit suffix to the metric header text used for
#     each data row's label (column A) and in the embedded charts' series
#     name cells. These are the same physical cells (A3, A4, A9, A10,
#     A14, A15, A19, A24) - updating the cell text here also updates every
#     chart series c:tx that references it.
$ws.Range("A3").Value = "encryption/$([char]0x3bc)s"
$ws.Range("A4").Value = "decryption/$([char]0x3bc)s"
$ws.Range("A9").Value = "encryption/$([char]0x3bc)s"
$ws.Range("A10").Value = "decryption/$([char]0x3bc)s"
$ws.Range("A14").Value = "encryption/$([char]0x3bc)s"
$ws.Range("A15").Value = "decryption/$([char]0x3bc)s"
$ws.Range("A19").Value = "digest generation/$([char]0x3bc)s"
$ws.Range("A24").Value = "signature generation/$([char]0x3bc)s"

# The former "Present" shared string is reused (without an explicit unit
# suffix) for the trailing label cell at A27.
$ws.Range("A27").Value = "Present"

# --- 8. Column A is now wider, to fit the longer "signature
#     generation/us" header text. The engine quantizes ColumnWidth to
#     whole-pixel (1/7 character-unit) steps on export, so an input of
#     21.6 is the closest achievable value to the saved file's recorded
#     width of 22.33203125 characters.
$ws.Columns("A").ColumnWidth = 21.6

# --- 9. Selection moved from U27 to I18 before the file was last saved.
$ws.Range("I18").Select()
